# Finished event table notes
# Fill in the Notes column (I) for rows 8-22 of the eventTable sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("eventTable")

# Rows whose Notes cell becomes "N/A" - these also need the same
# center-aligned / bold style already used by the other "N/A" cells
# in that row (column H has that style), so copy the format from H.
$naRows = @(8, 9, 11, 13, 15, 18, 19, 20, 21)
foreach ($r in $naRows) {
    $src = $ws.Cells.Item($r, 8)   # column H already styled as "N/A"
    $dst = $ws.Cells.Item($r, 9)   # column I (Notes)
    $src.Copy()
    $dst.PasteSpecial(-4122)       # xlPasteFormats
    $dst.Value = "N/A"
}
$excel.CutCopyMode = $false

# Rows with real note text (style stays as-is, just fill the value).
# Assigned in this order so new shared-string entries are appended in
# the same sequence as in the target workbook.
$ws.Cells.Item(14, 9).Value = "Used for withholding security deposits in the event that a customer returns a pack in a damaged or incomplete state"
$ws.Cells.Item(12, 9).Value = "Catches the possibility of a declined payment"
$ws.Cells.Item(10, 9).Value = "Charges the customer when they have requested to hire an equipment pack"
$ws.Cells.Item(16, 9).Value = "This event is for specifically marking that a pack has been collected by a customer after they have paid their hire fee and deposit"
$ws.Cells.Item(17, 9).Value = "Clearly marks the pack as being returned to a storage location"
$ws.Cells.Item(22, 9).Value = "This event enables the bank to send relevant transaction information to the ATO (see Rich Picture)"

# Update the selected cell shown in the sheet view.
$ws.Range("J26").Select()
